$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sectors")
$ws2 = $wb.Worksheets.Item("Fuels")

# --- Sectors sheet: reclassify fugitive / other-transformation sectors from
# Energy_Combustion (kt) to GDP (B2005USD), i.e. treat them as process emissions.
$ws1.Range("B5").Value = "GDP"
$ws1.Range("C5").Value = "B2005USD"

$ws1.Range("B30").Value = "GDP"
$ws1.Range("C30").Value = "B2005USD"

$ws1.Range("B31").Value = "GDP"
$ws1.Range("C31").Value = "B2005USD"

$ws1.Range("B32").Value = "GDP"
$ws1.Range("C32").Value = "B2005USD"

# --- Remove the natural-emission sector rows (11A/11B/11C) by clearing them.
$ws1.Range("A57:C59").ClearContents()

# --- Column widths for the Sectors sheet.
$ws1.Columns.Item(1).ColumnWidth = 28.33
$ws1.Columns.Item(2).ColumnWidth = 27.5

# --- Page setup: portrait orientation.
$ws1.PageSetup.Orientation = 1

# --- View state: Sectors becomes the active sheet, scrolled near the bottom
# with A56 selected.
$ws1.Activate()
$ws1.Range("A56").Select()
